{"js": "// Update candidate personal/contract data fields (status_impresion refresh\n// for candidate2): name, exterior number, city, daily salary, hire date and\n// e-mail throughout the generated contract document.\n\nasync function replaceAll(context, searchText, replacementText, options) {\n  const opts = Object.assign({ matchCase: true, matchWildcards: false }, options || {});\n  const results = context.document.body.search(searchText, opts);\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacementText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 1. Candidate full name.\nawait replaceAll(context, \"Nombre Completo: Carlos Fontaner\", \"Nombre Completo: Juan Alvarez\");\n\n// 2. Exterior number (remove the stray \"2\" suffix).\nawait replaceAll(context, \"N\u00famero Exterior: COLONIA2\", \"N\u00famero Exterior: COLONIA\");\n\n// 3. City field (label \"Ciudad:\" only \u2014 must not touch \"Estado: PUEBLA\").\nawait replaceAll(context, \"Ciudad: PUEBLA\", \"Ciudad: PUEBLAYORK\");\n\n// 4. Daily salary, every occurrence in the document (label line + clause text).\nawait replaceAll(context, \"$10000.00\", \"$5000.00\");\n\n// 5. Hire date, every occurrence (label line + contract-duration clause).\nawait replaceAll(context, \"2024-11-03\", \"2024-11-04\");\n\n// 6. Candidate e-mail address.\nawait replaceAll(context, \"Correo Electr\u00f3nico: carlitos@gmail.com\", \"Correo Electr\u00f3nico: juanito@gmail.com\");\n", "ps1": "# Update candidate personal/contract data fields (status_impresion refresh\n# for candidate2): name, exterior number, city, daily salary, hire date and\n# e-mail throughout the generated contract document.\n\n$d = $word.ActiveDocument\n\nfunction Replace-All($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# 1. Candidate full name.\nReplace-All 'Nombre Completo: Carlos Fontaner' 'Nombre Completo: Juan Alvarez'\n\n# 2. Exterior number (remove the stray \"2\" suffix).\nReplace-All 'N\u00famero Exterior: COLONIA2' 'N\u00famero Exterior: COLONIA'\n\n# 3. City field (label \"Ciudad:\" only \u2014 must not touch \"Estado: PUEBLA\").\nReplace-All 'Ciudad: PUEBLA' 'Ciudad: PUEBLAYORK'\n\n# 4. Daily salary, every occurrence in the document (label line + clause text).\nReplace-All '$10000.00' '$5000.00'\n\n# 5. Hire date, every occurrence (label line + contract-duration clause).\nReplace-All '2024-11-03' '2024-11-04'\n\n# 6. Candidate e-mail address.\nReplace-All 'Correo Electr\u00f3nico: carlitos@gmail.com' 'Correo Electr\u00f3nico: juanito@gmail.com'\n"}
